# Insert a new price-record row at row 45 (pushing the existing rows 45-52
# down to 46-53) and populate it with the new weekly Damasco record for
# "Vega Modelo de Temuco".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 45:52 down to 46:53, creating a blank row 45.
$ws.Rows("45:45").Insert()

# Populate the new row 45 with the new record.
$ws.Range("A45").Value = 10
$ws.Range("B45").Value = "Vega Modelo de Temuco"
$ws.Range("C45").Value = "La Araucanía"
$ws.Range("D45").Value = 44559
$ws.Range("E45").Value = 9
$ws.Range("F45").Value = "Fruta"
$ws.Range("G45").Value = 100103
$ws.Range("H45").Value = "Frutos de hueso (carozo)"
$ws.Range("I45").Value = 100103003
$ws.Range("J45").Value = "Damasco"
$ws.Range("K45").Value = "Modesto"
$ws.Range("L45").Value = "Primera"
$ws.Range("M45").Value = 95
$ws.Range("N45").Value = 18000
$ws.Range("O45").Value = 18000
$ws.Range("P45").Value = 18000
$ws.Range("Q45").Value = "$/bandeja 18 kilos"
$ws.Range("R45").Value = "Provincia de Quillota"
$ws.Range("S45").Value = 1000
$ws.Range("T45").Value = 18
